$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (46060 -> 46061) for every data row (rows 2 through 334).
$ws.Range("C2:C334").Value = 46061
